# Update the "Förändrad" (changed) date column (C) for all data rows
# (rows 2 through 116) from 2023-09-03 (serial 45172) to 2023-09-06
# (serial 45175).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C116").Value = 45175
